# Update the "data_validation" sheet: rename referenced source file
# from Contact_info.csv to Contact_info_21092024.csv for the CONTACT_INFO
# test rows (C2:C11), and move the active/top-left selection in the frozen
# pane back towards the start of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data_validation")

$ws.Range("C2:C11").Value = "Contact_info_21092024.csv"

$ws.Range("C6").Select()
